$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.180802702903748
$ws.Range("B1").Value = 2.419274806976318
$ws.Range("D1").Value = 2.32984471321106
$ws.Range("E1").Value = 1.192080497741699
